$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117, shifting existing rows 117:170 down to 118:171
$ws.Rows(117).Insert()

# Populate the newly inserted row 117 with the new data record
$ws.Range("A117").Value = 4
$ws.Range("B117").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C117").Value = "Los Lagos"
$ws.Range("D117").Value = 44510
$ws.Range("E117").Value = 10
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100102
$ws.Range("H117").Value = "Cítricos"
$ws.Range("I117").Value = 100102006
$ws.Range("J117").Value = "Pomelo"
$ws.Range("K117").Value = "Start Ruby"
$ws.Range("L117").Value = "Primera"
$ws.Range("M117").Value = 40
$ws.Range("N117").Value = 11000
$ws.Range("O117").Value = 12000
$ws.Range("P117").Value = 11500
$ws.Range("Q117").Value = "$/caja 14 kilos empedrada"
$ws.Range("R117").Value = "Región de O'Higgins"
$ws.Range("S117").Value = 821
$ws.Range("T117").Value = 14
